# Auto-update stock values: 2025-12-09 07:55:33 UTC
# Adds a new trading-day column (2025-12-09) to every per-field sheet,
# copying the header cell's format so the new date cell matches the
# existing header style, and filling the new data column values.

$wb = $excel.ActiveWorkbook

# --- Sheet "시가": add column BR (new date 20251209) ---
$ws = $wb.Worksheets.Item("시가")
$ws.Range("BQ1").Copy()
$ws.Range("BR1").PasteSpecial(-4122)
$ws.Columns.Item(70).ColumnWidth = 11.140625
$ws.Cells.Item(1,70).Value = 20251209
$ws.Cells.Item(2,70).Value = 15205
$ws.Cells.Item(3,70).Value = 57845
$ws.Cells.Item(4,70).Value = 18185
$ws.Cells.Item(5,70).Value = 16265
$ws.Cells.Item(6,70).Value = 17110
$ws.Cells.Item(7,70).Value = 19790
$ws.Cells.Item(8,70).Value = 18755
$ws.Cells.Item(9,70).Value = 1618
$ws.Cells.Item(10,70).Value = 17340
$ws.Cells.Item(11,70).Value = 5635
$ws.Cells.Item(12,70).Value = 6105
$ws.Cells.Item(13,70).Value = 6205
$ws.Cells.Item(14,70).Value = 1780
$ws.Cells.Item(15,70).Value = 19145
$ws.Cells.Item(16,70).Value = 11450
$ws.Cells.Item(17,70).Value = 6805
$ws.Cells.Item(18,70).Value = 18250
$ws.Cells.Item(19,70).Value = 16255

# --- Sheet "고가": add column BR (new date 20251209) ---
$ws = $wb.Worksheets.Item("고가")
$ws.Range("BQ1").Copy()
$ws.Range("BR1").PasteSpecial(-4122)
$ws.Columns.Item(70).ColumnWidth = 11.140625
$ws.Cells.Item(1,70).Value = 20251209
$ws.Cells.Item(2,70).Value = 15215
$ws.Cells.Item(3,70).Value = 58200
$ws.Cells.Item(4,70).Value = 18240
$ws.Cells.Item(5,70).Value = 16275
$ws.Cells.Item(6,70).Value = 17260
$ws.Cells.Item(7,70).Value = 19955
$ws.Cells.Item(8,70).Value = 18795
$ws.Cells.Item(9,70).Value = 1673
$ws.Cells.Item(10,70).Value = 17425
$ws.Cells.Item(11,70).Value = 5850
$ws.Cells.Item(12,70).Value = 6370
$ws.Cells.Item(13,70).Value = 6465
$ws.Cells.Item(14,70).Value = 1854
$ws.Cells.Item(15,70).Value = 19175
$ws.Cells.Item(16,70).Value = 11775
$ws.Cells.Item(17,70).Value = 7010
$ws.Cells.Item(18,70).Value = 18535
$ws.Cells.Item(19,70).Value = 16255

# --- Sheet "저가": add column BR (new date 20251209) ---
$ws = $wb.Worksheets.Item("저가")
$ws.Range("BQ1").Copy()
$ws.Range("BR1").PasteSpecial(-4122)
$ws.Columns.Item(70).ColumnWidth = 11.140625
$ws.Cells.Item(1,70).Value = 20251209
$ws.Cells.Item(2,70).Value = 14900
$ws.Cells.Item(3,70).Value = 57405
$ws.Cells.Item(4,70).Value = 17885
$ws.Cells.Item(5,70).Value = 16045
$ws.Cells.Item(6,70).Value = 17060
$ws.Cells.Item(7,70).Value = 19645
$ws.Cells.Item(8,70).Value = 18515
$ws.Cells.Item(9,70).Value = 1607
$ws.Cells.Item(10,70).Value = 17185
$ws.Cells.Item(11,70).Value = 5620
$ws.Cells.Item(12,70).Value = 6090
$ws.Cells.Item(13,70).Value = 6150
$ws.Cells.Item(14,70).Value = 1763
$ws.Cells.Item(15,70).Value = 18800
$ws.Cells.Item(16,70).Value = 11380
$ws.Cells.Item(17,70).Value = 6770
$ws.Cells.Item(18,70).Value = 18250
$ws.Cells.Item(19,70).Value = 15970

# --- Sheet "종가": add column BR (new date 20251209) ---
$ws = $wb.Worksheets.Item("종가")
$ws.Range("BQ1").Copy()
$ws.Range("BR1").PasteSpecial(-4122)
$ws.Columns.Item(70).ColumnWidth = 11.140625
$ws.Cells.Item(1,70).Value = 20251209
$ws.Cells.Item(2,70).Value = 15155
$ws.Cells.Item(3,70).Value = 57600
$ws.Cells.Item(4,70).Value = 18120
$ws.Cells.Item(5,70).Value = 16115
$ws.Cells.Item(6,70).Value = 17145
$ws.Cells.Item(7,70).Value = 19790
$ws.Cells.Item(8,70).Value = 18610
$ws.Cells.Item(9,70).Value = 1637
$ws.Cells.Item(10,70).Value = 17400
$ws.Cells.Item(11,70).Value = 5785
$ws.Cells.Item(12,70).Value = 6290
$ws.Cells.Item(13,70).Value = 6360
$ws.Cells.Item(14,70).Value = 1814
$ws.Cells.Item(15,70).Value = 19050
$ws.Cells.Item(16,70).Value = 11645
$ws.Cells.Item(17,70).Value = 6935
$ws.Cells.Item(18,70).Value = 18415
$ws.Cells.Item(19,70).Value = 16040

# --- Sheet "거래량": add column BR (new date 20251209) ---
$ws = $wb.Worksheets.Item("거래량")
$ws.Range("BQ1").Copy()
$ws.Range("BR1").PasteSpecial(-4122)
$ws.Columns.Item(70).ColumnWidth = 11.140625
$ws.Cells.Item(1,70).Value = 20251209
$ws.Cells.Item(2,70).Value = 191632
$ws.Cells.Item(3,70).Value = 1203651
$ws.Cells.Item(4,70).Value = 122375
$ws.Cells.Item(5,70).Value = 113962
$ws.Cells.Item(6,70).Value = 1117997
$ws.Cells.Item(7,70).Value = 2845145
$ws.Cells.Item(8,70).Value = 2675787
$ws.Cells.Item(9,70).Value = 5439163
$ws.Cells.Item(10,70).Value = 321572
$ws.Cells.Item(11,70).Value = 1168004
$ws.Cells.Item(12,70).Value = 4429980
$ws.Cells.Item(13,70).Value = 11622331
$ws.Cells.Item(14,70).Value = 43206720
$ws.Cells.Item(15,70).Value = 1361571
$ws.Cells.Item(16,70).Value = 569720
$ws.Cells.Item(17,70).Value = 76967
$ws.Cells.Item(18,70).Value = 61685
$ws.Cells.Item(19,70).Value = 54276

# --- Sheet "s20": add column AY (new date 20251209) ---
$ws = $wb.Worksheets.Item("s20")
$ws.Range("AX1").Copy()
$ws.Range("AY1").PasteSpecial(-4122)
$ws.Columns.Item(51).ColumnWidth = 9.140625
$ws.Cells.Item(1,51).Value = 20251209
$ws.Cells.Item(2,51).Value = 11
$ws.Cells.Item(3,51).Value = 71
$ws.Cells.Item(4,51).Value = 11
$ws.Cells.Item(5,51).Value = 61
$ws.Cells.Item(6,51).Value = 90
$ws.Cells.Item(7,51).Value = 98
$ws.Cells.Item(8,51).Value = 74
$ws.Cells.Item(9,51).Value = 56
$ws.Cells.Item(10,51).Value = 84
$ws.Cells.Item(11,51).Value = 100
$ws.Cells.Item(12,51).Value = 100
$ws.Cells.Item(13,51).Value = 100
$ws.Cells.Item(14,51).Value = 70
$ws.Cells.Item(15,51).Value = 39
$ws.Cells.Item(16,51).Value = 99
$ws.Cells.Item(17,51).Value = 89
$ws.Cells.Item(18,51).Value = 84
$ws.Cells.Item(19,51).Value = 55

# --- Sheet "s60": add column K (new date 20251209) ---
$ws = $wb.Worksheets.Item("s60")
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Columns.Item(11).ColumnWidth = 9.140625
$ws.Cells.Item(1,11).Value = 20251209
$ws.Cells.Item(2,11).Value = 64
$ws.Cells.Item(3,11).Value = 83
$ws.Cells.Item(4,11).Value = 69
$ws.Cells.Item(5,11).Value = 78
$ws.Cells.Item(6,11).Value = 91
$ws.Cells.Item(7,11).Value = 93
$ws.Cells.Item(8,11).Value = 88
$ws.Cells.Item(9,11).Value = 64
$ws.Cells.Item(10,11).Value = 85
$ws.Cells.Item(11,11).Value = 100
$ws.Cells.Item(12,11).Value = 100
$ws.Cells.Item(13,11).Value = 100
$ws.Cells.Item(14,11).Value = 74
$ws.Cells.Item(15,11).Value = 78
$ws.Cells.Item(16,11).Value = 92
$ws.Cells.Item(17,11).Value = 89
$ws.Cells.Item(18,11).Value = 95
$ws.Cells.Item(19,11).Value = 83

# --- Sheet "z20": add column AY (new date 20251209) ---
$ws = $wb.Worksheets.Item("z20")
$ws.Range("AX1").Copy()
$ws.Range("AY1").PasteSpecial(-4122)
$ws.Columns.Item(51).ColumnWidth = 9.140625
$ws.Cells.Item(1,51).Value = 20251209
$ws.Cells.Item(2,51).Value = -93
$ws.Cells.Item(3,51).Value = 31
$ws.Cells.Item(4,51).Value = -58
$ws.Cells.Item(5,51).Value = 11
$ws.Cells.Item(6,51).Value = 49
$ws.Cells.Item(7,51).Value = 65
$ws.Cells.Item(8,51).Value = 37
$ws.Cells.Item(9,51).Value = 32
$ws.Cells.Item(10,51).Value = 31
$ws.Cells.Item(11,51).Value = 99
$ws.Cells.Item(12,51).Value = 101
$ws.Cells.Item(13,51).Value = 108
$ws.Cells.Item(14,51).Value = 52
$ws.Cells.Item(15,51).Value = -18
$ws.Cells.Item(16,51).Value = 85
$ws.Cells.Item(17,51).Value = 77
$ws.Cells.Item(18,51).Value = 33
$ws.Cells.Item(19,51).Value = 8

# --- Sheet "z60": add column K (new date 20251209) ---
$ws = $wb.Worksheets.Item("z60")
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Columns.Item(11).ColumnWidth = 9.140625
$ws.Cells.Item(1,11).Value = 20251209
$ws.Cells.Item(2,11).Value = 34
$ws.Cells.Item(3,11).Value = 39
$ws.Cells.Item(4,11).Value = 38
$ws.Cells.Item(5,11).Value = 26
$ws.Cells.Item(6,11).Value = 50
$ws.Cells.Item(7,11).Value = 56
$ws.Cells.Item(8,11).Value = 43
$ws.Cells.Item(9,11).Value = 36
$ws.Cells.Item(10,11).Value = 43
$ws.Cells.Item(11,11).Value = 71
$ws.Cells.Item(12,11).Value = 70
$ws.Cells.Item(13,11).Value = 78
$ws.Cells.Item(14,11).Value = 47
$ws.Cells.Item(15,11).Value = 60
$ws.Cells.Item(16,11).Value = 60
$ws.Cells.Item(17,11).Value = 58
$ws.Cells.Item(18,11).Value = 58
$ws.Cells.Item(19,11).Value = 68

# --- Sheet "gap": add column AY (new date 20251209) ---
$ws = $wb.Worksheets.Item("gap")
$ws.Range("AX1").Copy()
$ws.Range("AY1").PasteSpecial(-4122)
$ws.Columns.Item(51).ColumnWidth = 11.140625
$ws.Cells.Item(1,51).Value = "20251209"
$ws.Cells.Item(2,51).Value = 95
$ws.Cells.Item(3,51).Value = 102
$ws.Cells.Item(4,51).Value = 97
$ws.Cells.Item(5,51).Value = 101
$ws.Cells.Item(6,51).Value = 103
$ws.Cells.Item(7,51).Value = 105
$ws.Cells.Item(8,51).Value = 103
$ws.Cells.Item(9,51).Value = 105
$ws.Cells.Item(10,51).Value = 102
$ws.Cells.Item(11,51).Value = 110
$ws.Cells.Item(12,51).Value = 111
$ws.Cells.Item(13,51).Value = 114
$ws.Cells.Item(14,51).Value = 109
$ws.Cells.Item(15,51).Value = 99
$ws.Cells.Item(16,51).Value = 108
$ws.Cells.Item(17,51).Value = 108
$ws.Cells.Item(18,51).Value = 102
$ws.Cells.Item(19,51).Value = 100

# --- Sheet "std": add column AF (new date 20251209) ---
$ws = $wb.Worksheets.Item("std")
$ws.Range("AE1").Copy()
$ws.Range("AF1").PasteSpecial(-4122)
$ws.Columns.Item(32).ColumnWidth = 11.140625
$ws.Cells.Item(1,32).Value = "20251209"
$ws.Cells.Item(2,32).Value = -46.6
$ws.Cells.Item(3,32).Value = -7.09
$ws.Cells.Item(4,32).Value = -43.96
$ws.Cells.Item(5,32).Value = -7.96
$ws.Cells.Item(6,32).Value = -10.4
$ws.Cells.Item(7,32).Value = -2.41
$ws.Cells.Item(8,32).Value = -4.54
$ws.Cells.Item(9,32).Value = -8.18
$ws.Cells.Item(10,32).Value = -3.71
$ws.Cells.Item(11,32).Value = 19.82
$ws.Cells.Item(12,32).Value = 18.41
$ws.Cells.Item(13,32).Value = 29.93
$ws.Cells.Item(14,32).Value = -6.09
$ws.Cells.Item(15,32).Value = -48.79
$ws.Cells.Item(16,32).Value = 2.83
$ws.Cells.Item(17,32).Value = 5.47
$ws.Cells.Item(18,32).Value = -3.13
$ws.Cells.Item(19,32).Value = -49.11

# --- Sheet "quant": add column K (new date 20251209) ---
$ws = $wb.Worksheets.Item("quant")
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Columns.Item(11).ColumnWidth = 11.140625
$ws.Cells.Item(1,11).Value = "20251209"
$ws.Cells.Item(2,11).Value = 90
$ws.Cells.Item(3,11).Value = 26
$ws.Cells.Item(4,11).Value = 40
$ws.Cells.Item(5,11).Value = 7
$ws.Cells.Item(6,11).Value = 28
$ws.Cells.Item(7,11).Value = 37
$ws.Cells.Item(8,11).Value = 19
$ws.Cells.Item(9,11).Value = 38
$ws.Cells.Item(10,11).Value = 11
$ws.Cells.Item(11,11).Value = 47
$ws.Cells.Item(12,11).Value = 81
$ws.Cells.Item(13,11).Value = 89
$ws.Cells.Item(14,11).Value = 43
$ws.Cells.Item(15,11).Value = 42
$ws.Cells.Item(16,11).Value = 47
$ws.Cells.Item(17,11).Value = 16
$ws.Cells.Item(18,11).Value = 17
$ws.Cells.Item(19,11).Value = 57
